# Update the division problems in the single table of the worksheet.
# Each problem lives in its own table cell, and every old value is
# unique within the document, but a couple of the *new* values happen
# to equal *other* old values elsewhere in the table (e.g. "17÷6=" is
# both an old value in row 1 and a new value written into row 5, and
# "30÷4=" is both an old value in row 3 and a new value written into
# row 9). To guarantee each replacement only touches the intended
# cell (and never accidentally matches text written by an earlier
# step, nor gets re-matched by a later step), we scope every
# Find/Replace to that specific cell's Range rather than searching
# the whole document.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of (row, col) -> (old text, new text), taken from the diff.
$changes = @(
    @{ Row = 1;  Col = 1; Old = "17÷6="; New = "46÷3=" },
    @{ Row = 1;  Col = 2; Old = "40÷9="; New = "68÷7=" },
    @{ Row = 1;  Col = 3; Old = "49÷8="; New = "86÷8=" },
    @{ Row = 1;  Col = 4; Old = "30÷4="; New = "43÷8=" },
    @{ Row = 1;  Col = 5; Old = "61÷5="; New = "19÷5=" },

    @{ Row = 5;  Col = 1; Old = "75÷4="; New = "17÷6=" },
    @{ Row = 5;  Col = 2; Old = "99÷5="; New = "21÷6=" },
    @{ Row = 5;  Col = 3; Old = "11÷6="; New = "67÷9=" },
    @{ Row = 5;  Col = 4; Old = "72÷4="; New = "88÷2=" },
    @{ Row = 5;  Col = 5; Old = "13÷4="; New = "98÷4=" },

    @{ Row = 9;  Col = 1; Old = "32÷2="; New = "26÷4=" },
    @{ Row = 9;  Col = 2; Old = "99÷9="; New = "80÷8=" },
    @{ Row = 9;  Col = 3; Old = "28÷2="; New = "10÷3=" },
    @{ Row = 9;  Col = 4; Old = "53÷5="; New = "30÷4=" },
    @{ Row = 9;  Col = 5; Old = "20÷7="; New = "23÷8=" },

    @{ Row = 13; Col = 1; Old = "70÷5="; New = "39÷7=" },
    @{ Row = 13; Col = 2; Old = "45÷3="; New = "13÷9=" },
    @{ Row = 13; Col = 3; Old = "73÷4="; New = "14÷5=" },
    @{ Row = 13; Col = 4; Old = "96÷9="; New = "81÷3=" },
    @{ Row = 13; Col = 5; Old = "37÷5="; New = "46÷4=" },

    @{ Row = 17; Col = 1; Old = "39÷6="; New = "14÷8=" },
    @{ Row = 17; Col = 2; Old = "64÷6="; New = "47÷9=" },
    @{ Row = 17; Col = 3; Old = "84÷4="; New = "61÷2=" },
    @{ Row = 17; Col = 4; Old = "55÷6="; New = "63÷5=" },
    @{ Row = 17; Col = 5; Old = "99÷6="; New = "21÷7=" }
)

foreach ($change in $changes) {
    $cell = $t.Cell($change.Row, $change.Col)
    $rng = $cell.Range
    $rng.Find.Execute($change.Old, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $change.New, 2) | Out-Null
}
